# Pandas -- Updated the dummy class data files.
# Removed dummy_class and replaced it with two files with the data split.
#
# This workbook is "file 2 of 2": the leading descriptive text row is
# removed, and the student rows that have no "Test 3 (50%)" score
# (i.e. those that were split out into the other file) are deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from the bottom up so earlier row numbers stay valid:
#   row 1  -> the descriptive "This is a dummy set of data..." text
#   row 20 -> A7210476B (no Test 3 score)
#   row 31 -> A3699958T (no Test 3 score)
#   row 35 -> A6867791C (no Test 3 score)
#   row 37 -> A7667457P (no Test 3 score)
$rowsToDelete = @(37, 35, 31, 20, 1)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).EntireRow.Delete()
}

# Touch a few trailing rows far below the data (mirrors the source
# workbook's saved state, which tracks a used range extending to the
# bottom of the sheet with a handful of slightly-shorter rows there).
$tail = $ws.Range("A1048572:B1048576")
$tail.Value = "x"
$tail.ClearContents()
$tail.RowHeight = 12.8

# Restore the active selection.
$ws.Range("U18").Select()
